{"js": "// Oficios Primer Parcial 20221\n// Update the \"SUPERVISA EL CUMPLIMIENTO DE TAREAS Y PROCESOS PARA EVALUAR LA\n// PRODUCTIVIDAD EN LA ORGANIZACI\u00d3N\" row: reprobados 10 -> 8 (and its\n// percentage 25.0% -> 20.0%), the matching asesor\u00eda count 10 -> 8, and the\n// \"N\u00famero de Alumnos\" summary counts 33 -> 32 / 1 -> 2.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Table index 1: ASIGNATURA / NO. DE ALUMNOS REPROBADOS / PORCENTAJE DE REPROBADOS\nconst reprobadosTable = tables.items[1];\nreprobadosTable.getCell(5, 1).value = \"8\";\nreprobadosTable.getCell(5, 2).value = \"20.0%\";\n\n// Table index 2: ASIGNATURA DE ASESOR\u00cdA / NO. DE ALUMNOS CON ASESOR\u00cdA / NOMBRE DEL ASESOR\nconst asesoriaTable = tables.items[2];\nasesoriaTable.getCell(5, 1).value = \"8\";\n\n// Table index 3: N\u00famero de Alumnos que tienen asesor\u00eda / canalizados / no requirieron\nconst resumenTable = tables.items[3];\nresumenTable.getCell(0, 4).value = \"32\"; // N\u00famero de Alumnos que tienen asesor\u00eda (M)\nresumenTable.getCell(2, 4).value = \"2\";  // N\u00famero de Alumnos que no requirieron atenci\u00f3n del tutor (M)\n\nawait context.sync();\n", "ps1": "# Oficios Primer Parcial 20221\n# Update the \"SUPERVISA EL CUMPLIMIENTO DE TAREAS Y PROCESOS PARA EVALUAR LA\n# PRODUCTIVIDAD EN LA ORGANIZACI\u00d3N\" row: reprobados 10 -> 8 (and its\n# percentage 25.0% -> 20.0%), the matching asesor\u00eda count 10 -> 8, and the\n# \"N\u00famero de Alumnos\" summary counts 33 -> 32 / 1 -> 2.\n\n$d = $word.ActiveDocument\n\n# Table 2 (1-indexed): ASIGNATURA / NO. DE ALUMNOS REPROBADOS / PORCENTAJE DE REPROBADOS\n$reprobadosTable = $d.Tables.Item(2)\n$reprobadosTable.Cell(6, 2).Range.Text = \"8\"\n$reprobadosTable.Cell(6, 3).Range.Text = \"20.0%\"\n\n# Table 3 (1-indexed): ASIGNATURA DE ASESOR\u00cdA / NO. DE ALUMNOS CON ASESOR\u00cdA / NOMBRE DEL ASESOR\n$asesoriaTable = $d.Tables.Item(3)\n$asesoriaTable.Cell(6, 2).Range.Text = \"8\"\n\n# Table 4 (1-indexed): N\u00famero de Alumnos que tienen asesor\u00eda / canalizados / no requirieron\n$resumenTable = $d.Tables.Item(4)\n$resumenTable.Cell(1, 5).Range.Text = \"32\"\n$resumenTable.Cell(3, 5).Range.Text = \"2\"\n"}
